# Apply updated TPM values for the Fgf1-Fgfr2 LR-pair sheet.
# The cluster set changed (Inflammatory-Mac dropped, Resolving-Mac added as a
# sending cluster only) and every numeric column was recomputed, shrinking the
# data block from 16 rows (4 senders x 4 targets) down to 12 (4 senders x 3 targets).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,20
$row2[0,0] = 'ECs'
$row2[0,1] = 'Fgf1'
$row2[0,2] = 'Fgfr2'
$row2[0,3] = 'ECs'
$row2[0,4] = 2
$row2[0,5] = 0.6666666666666666
$row2[0,6] = 0.265349
$row2[0,7] = 0.7960469999999999
$row2[0,8] = 0.1498685997319469
$row2[0,9] = 0.1498685997319469
$row2[0,10] = 3
$row2[0,11] = 1
$row2[0,12] = 0.8155003333333334
$row2[0,13] = 2.446501
$row2[0,14] = 0.1910612426590028
$row2[0,15] = 0.1910612426590029
$row2[0,16] = 0.2163921979496667
$row2[0,17] = 1.947529781547
$row2[0,18] = 0.02863408090035047
$row2[0,19] = 0.02863408090035047
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object 'object[,]' 1,20
$row3[0,0] = 'ECs'
$row3[0,1] = 'Fgf1'
$row3[0,2] = 'Fgfr2'
$row3[0,3] = 'FAPs'
$row3[0,4] = 2
$row3[0,5] = 0.6666666666666666
$row3[0,6] = 0.265349
$row3[0,7] = 0.7960469999999999
$row3[0,8] = 0.1498685997319469
$row3[0,9] = 0.1498685997319469
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 3.333134333333334
$row3[0,13] = 9.999403000000001
$row3[0,14] = 0.7809105179307759
$row3[0,15] = 0.780910517930776
$row3[0,16] = 0.8844438622156667
$row3[0,17] = 7.959994759941
$row3[0,18] = 0.1170339658382348
$row3[0,19] = 0.1170339658382348
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object 'object[,]' 1,20
$row4[0,0] = 'ECs'
$row4[0,1] = 'Fgf1'
$row4[0,2] = 'Fgfr2'
$row4[0,3] = 'MuSCs'
$row4[0,4] = 2
$row4[0,5] = 0.6666666666666666
$row4[0,6] = 0.265349
$row4[0,7] = 0.7960469999999999
$row4[0,8] = 0.1498685997319469
$row4[0,9] = 0.1498685997319469
$row4[0,10] = 3
$row4[0,11] = 1
$row4[0,12] = 0.119632
$row4[0,13] = 0.358896
$row4[0,14] = 0.02802823941022116
$row4[0,15] = 0.02802823941022117
$row4[0,16] = 0.031744231568
$row4[0,17] = 0.285698084112
$row4[0,18] = 0.004200552993361613
$row4[0,19] = 0.004200552993361614
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object 'object[,]' 1,20
$row5[0,0] = 'FAPs'
$row5[0,1] = 'Fgf1'
$row5[0,2] = 'Fgfr2'
$row5[0,3] = 'ECs'
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 1.003400666666667
$row5[0,7] = 3.010202
$row5[0,8] = 0.5667187473230926
$row5[0,9] = 0.5667187473230925
$row5[0,10] = 3
$row5[0,11] = 1
$row5[0,12] = 0.8155003333333334
$row5[0,13] = 2.446501
$row5[0,14] = 0.1910612426590028
$row5[0,15] = 0.1910612426590029
$row5[0,16] = 0.8182735781335555
$row5[0,17] = 7.364462203201999
$row5[0,18] = 0.1082779881017035
$row5[0,19] = 0.1082779881017035
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object 'object[,]' 1,20
$row6[0,0] = 'FAPs'
$row6[0,1] = 'Fgf1'
$row6[0,2] = 'Fgfr2'
$row6[0,3] = 'FAPs'
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 1.003400666666667
$row6[0,7] = 3.010202
$row6[0,8] = 0.5667187473230926
$row6[0,9] = 0.5667187473230925
$row6[0,10] = 3
$row6[0,11] = 1
$row6[0,12] = 3.333134333333334
$row6[0,13] = 9.999403000000001
$row6[0,14] = 0.7809105179307759
$row6[0,15] = 0.780910517930776
$row6[0,16] = 3.344469212156223
$row6[0,17] = 30.100222909406
$row6[0,18] = 0.4425566304931567
$row6[0,19] = 0.4425566304931567
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object 'object[,]' 1,20
$row7[0,0] = 'FAPs'
$row7[0,1] = 'Fgf1'
$row7[0,2] = 'Fgfr2'
$row7[0,3] = 'MuSCs'
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 1.003400666666667
$row7[0,7] = 3.010202
$row7[0,8] = 0.5667187473230926
$row7[0,9] = 0.5667187473230925
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 0.119632
$row7[0,13] = 0.358896
$row7[0,14] = 0.02802823941022116
$row7[0,15] = 0.02802823941022117
$row7[0,16] = 0.1200388285546667
$row7[0,17] = 1.080349456992
$row7[0,18] = 0.01588412872823227
$row7[0,19] = 0.01588412872823227
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object 'object[,]' 1,20
$row8[0,0] = 'MuSCs'
$row8[0,1] = 'Fgf1'
$row8[0,2] = 'Fgfr2'
$row8[0,3] = 'ECs'
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 0.4790693333333333
$row8[0,7] = 1.437208
$row8[0,8] = 0.2705774288246194
$row8[0,9] = 0.2705774288246194
$row8[0,10] = 3
$row8[0,11] = 1
$row8[0,12] = 0.8155003333333334
$row8[0,13] = 2.446501
$row8[0,14] = 0.1910612426590028
$row8[0,15] = 0.1910612426590029
$row8[0,16] = 0.3906812010231111
$row8[0,17] = 3.516130809208
$row8[0,18] = 0.05169685978670969
$row8[0,19] = 0.0516968597867097
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object 'object[,]' 1,20
$row9[0,0] = 'MuSCs'
$row9[0,1] = 'Fgf1'
$row9[0,2] = 'Fgfr2'
$row9[0,3] = 'FAPs'
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 0.4790693333333333
$row9[0,7] = 1.437208
$row9[0,8] = 0.2705774288246194
$row9[0,9] = 0.2705774288246194
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 3.333134333333334
$row9[0,13] = 9.999403000000001
$row9[0,14] = 0.7809105179307759
$row9[0,15] = 0.780910517930776
$row9[0,16] = 1.596802442980444
$row9[0,17] = 14.371221986824
$row9[0,18] = 0.2112967600838112
$row9[0,19] = 0.2112967600838113
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object 'object[,]' 1,20
$row10[0,0] = 'MuSCs'
$row10[0,1] = 'Fgf1'
$row10[0,2] = 'Fgfr2'
$row10[0,3] = 'MuSCs'
$row10[0,4] = 3
$row10[0,5] = 1
$row10[0,6] = 0.4790693333333333
$row10[0,7] = 1.437208
$row10[0,8] = 0.2705774288246194
$row10[0,9] = 0.2705774288246194
$row10[0,10] = 3
$row10[0,11] = 1
$row10[0,12] = 0.119632
$row10[0,13] = 0.358896
$row10[0,14] = 0.02802823941022116
$row10[0,15] = 0.02802823941022117
$row10[0,16] = 0.05731202248533333
$row10[0,17] = 0.515808202368
$row10[0,18] = 0.00758380895409851
$row10[0,19] = 0.007583808954098511
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object 'object[,]' 1,20
$row11[0,0] = 'Resolving-Mac'
$row11[0,1] = 'Fgf1'
$row11[0,2] = 'Fgfr2'
$row11[0,3] = 'ECs'
$row11[0,4] = 1
$row11[0,5] = 0.3333333333333333
$row11[0,6] = 0.02272533333333333
$row11[0,7] = 0.068176
$row11[0,8] = 0.01283522412034115
$row11[0,9] = 0.01283522412034115
$row11[0,10] = 3
$row11[0,11] = 1
$row11[0,12] = 0.8155003333333334
$row11[0,13] = 2.446501
$row11[0,14] = 0.1910612426590028
$row11[0,15] = 0.1910612426590029
$row11[0,16] = 0.01853251690844445
$row11[0,17] = 0.166792652176
$row11[0,18] = 0.002452313870239186
$row11[0,19] = 0.002452313870239187
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object 'object[,]' 1,20
$row12[0,0] = 'Resolving-Mac'
$row12[0,1] = 'Fgf1'
$row12[0,2] = 'Fgfr2'
$row12[0,3] = 'FAPs'
$row12[0,4] = 1
$row12[0,5] = 0.3333333333333333
$row12[0,6] = 0.02272533333333333
$row12[0,7] = 0.068176
$row12[0,8] = 0.01283522412034115
$row12[0,9] = 0.01283522412034115
$row12[0,10] = 3
$row12[0,11] = 1
$row12[0,12] = 3.333134333333334
$row12[0,13] = 9.999403000000001
$row12[0,14] = 0.7809105179307759
$row12[0,15] = 0.780910517930776
$row12[0,16] = 0.07574658876977779
$row12[0,17] = 0.6817192989280001
$row12[0,18] = 0.01002316151557319
$row12[0,19] = 0.01002316151557319
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object 'object[,]' 1,20
$row13[0,0] = 'Resolving-Mac'
$row13[0,1] = 'Fgf1'
$row13[0,2] = 'Fgfr2'
$row13[0,3] = 'MuSCs'
$row13[0,4] = 1
$row13[0,5] = 0.3333333333333333
$row13[0,6] = 0.02272533333333333
$row13[0,7] = 0.068176
$row13[0,8] = 0.01283522412034115
$row13[0,9] = 0.01283522412034115
$row13[0,10] = 3
$row13[0,11] = 1
$row13[0,12] = 0.119632
$row13[0,13] = 0.358896
$row13[0,14] = 0.02802823941022116
$row13[0,15] = 0.02802823941022117
$row13[0,16] = 0.002718677077333333
$row13[0,17] = 0.024468093696
$row13[0,18] = 0.000359748734528767
$row13[0,19] = 0.000359748734528767
$ws.Range("A13:T13").Value = $row13

# The source table now has only 12 data rows (A2:T13); remove the 4 leftover
# rows from the old 16-row block so the sheet dimension matches (A1:T13).
$ws.Range("A14:T17").Delete() | Out-Null
